$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/centered table style (currently ending at row 12) down
# to the four newly-added data rows (13-16), matching the existing pattern
# used for rows 5-12.
$ws.Range("C12:F12").Copy()
$ws.Range("C13:F16").PasteSpecial(-4122)

# Extend the plain blank-row style (currently ending at row 26) down to the
# two newly-added blank rows (27-28).
$ws.Range("C26:F26").Copy()
$ws.Range("C27:F28").PasteSpecial(-4122)

# --- Row 5: "Implementacija ASP NET CORE Web Api + deploy" moved here from
#     row 12's old position, with updated dates/description.
$ws.Range("C5").Value = "Implementacija ASP NET CORE Web Api + deploy"
$ws.Range("D5").Value = "20.04.2024."
$ws.Range("E5").Value = "22.04.2024."
$ws.Range("F5").Value = "Endpointovi za pohranu podataka historiju podataka zabilježenih sa senzora i deploy web api-a i baze podataka na fitov server."

# --- Row 6: "Frontend implementacija"
$ws.Range("C6").Value = "Frontend implementacija"
$ws.Range("D6").Value = "22.04.2024."
$ws.Range("E6").Value = "24.04.2024."
$ws.Range("F6").Value = "Komponente temperature/humidity, security, gas/fire i korisnicki profil"

# --- Row 7: "Dodavanje koda u Arduino IDE"
$ws.Range("C7").Value = "Dodavanje koda u Arduino IDE"
$ws.Range("D7").Value = "25.04.2024."
$ws.Range("E7").Value = "25.04.2024."
$ws.Range("F7").Value = "Implementacija slanja JSON objekata na napravljenje endpointove i testiranje komunikacije na frontendu"

# --- Row 8: "Implementacija notifikacija email, whatsApp i sms poruke na web api-u"
$ws.Range("C8").Value = "Implementacija notifikacija email, whatsApp i sms poruke na web api-u"
$ws.Range("D8").Value = "26.04.2024."
$ws.Range("E8").Value = "27.04.2024."
$ws.Range("F8").Value = "Koristenje Smtp servisa i Twilio servisa u c# asp net core web api"

# --- Row 9: "Dodavanje senzora light"
$ws.Range("C9").Value = "Dodavanje senzora light"
$ws.Range("D9").Value = "27.04.2024."
$ws.Range("E9").Value = "28.04.2024."
$ws.Range("F9").Value = "Dodavanje koda u Arduino IDE za upravljanje light senzorom, izmjene u Realtime DB na Firebase, dodavanje Light komponente u Angularu i testiranje komunikacije na frontendu."

# --- Row 10: "Implementacija mogucnosti dodavanja vise uredjaja - AUTENTIFIKACIJA"
$ws.Range("C10").Value = "Implementacija mogucnosti dodavanja vise uredjaja - AUTENTIFIKACIJA"
$ws.Range("D10").Value = "04.05.2024."
$ws.Range("E10").Value = "04.05.2024."
$ws.Range("F10").Value = "Korekcija table u sql bazi podataka i real-time bazi podataka na firebase, kreiranje endpintova na Web API (autentifikacija)"

# --- Row 11: "Frontend implementacija - AUTENTIFIKACIJA"
$ws.Range("C11").Value = "Frontend implementacija - AUTENTIFIKACIJA"
$ws.Range("D11").Value = "05.05.2024."
$ws.Range("E11").Value = "05.05.2024."
$ws.Range("F11").Value = "Korekcija komponente Profil u Angular projektu, dodavanje AuthServisa, Login komponente i Password komponente za promjenu sifre"

# --- Row 12: "Izmjene u zaglavlju endpointova i firebase upita u Arduino IDE"
#     (F12 stays blank, as it already was.)
$ws.Range("C12").Value = "Izmjene u zaglavlju endpointova i firebase upita u Arduino IDE"
$ws.Range("D12").Value = "05.05.2024."
$ws.Range("E12").Value = "05.05.2024."

# --- Row 13 (new): "Dodavanje pregleda prijavljenih uredjaja i mogucnost odjave"
$ws.Range("C13").Value = "Dodavanje pregleda prijavljenih uredjaja i mogucnost odjave"
$ws.Range("D13").Value = "06.05.2024."
$ws.Range("E13").Value = "06.05.2024."
$ws.Range("F13").Value = "Dodavanje endpointova, dodavanje komponente u angularu"

# --- Row 14 (new): "Optimizacija koda na frontendu i u Arduino IDE"
$ws.Range("C14").Value = "Optimizacija koda na frontendu i u Arduino IDE"
$ws.Range("D14").Value = "07.05.2024."
$ws.Range("E14").Value = "13.05.2024."

# --- Row 15 (new): "Dodavanje dark/light mode"
$ws.Range("C15").Value = "Dodavanje dark/light mode"
$ws.Range("D15").Value = "13.05.2024."
$ws.Range("E15").Value = "14.05.2024."
$ws.Range("F15").Value = "Dodavanje photoresistora na nodemcu, dodavanje varijable na firebase i izmjene na frontendu"

# --- Row 16 (new): "Deploy  web stranice na github.io i app.fit.ba "
$ws.Range("C16").Value = "Deploy  web stranice na github.io i app.fit.ba "
$ws.Range("D16").Value = "14.05.2024."
$ws.Range("E16").Value = "14.05.2024."

# Update the selected cell to match the saved view state.
$ws.Range("C10").Select()
